# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status for fe2edfbc row changes from "Ready for handoff" to "Handback transform failed"
# (same shared string is used on the Overview summary sheet and both detail sheets)
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# New error detail messages
$zhCnError = "Handback file name: qu2mygtu.s3q is different with handoff file name: fe2edfbc-f63f-4a1f-b9fa-1c465fdfdd36.d1c03207ca8ae36ed0ee2c97f5edf1d65e624abb.zh-cn."
$deDeError = "Handback file name: qu2mygtu.s3q is different with handoff file name: fe2edfbc-f63f-4a1f-b9fa-1c465fdfdd36.d1c03207ca8ae36ed0ee2c97f5edf1d65e624abb.de-de."

$wsZhCn.Range("P3").Value = $zhCnError
$wsDeDe.Range("P3").Value = $deDeError

# Widen the "Error Detail" column (P) on both sheets so the message is readable.
# NOTE: the host's ColumnWidth setter re-derives the stored OOXML width from a
# pixel-quantised character width, so asking for exactly 40 round-trips to
# 40.8333. 39.16 is the input that round-trips to a clean stored width of 40.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.16
$wsDeDe.Columns.Item(16).ColumnWidth = 39.16
